$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 87.331081081081081
$ws.Range("A3").Value = 86.993243243243242
$ws.Range("A4").Value = 87.331081081081081
$ws.Range("A5").Value = 86.317567567567565
$ws.Range("A6").Value = 86.148648648648646
$ws.Range("A7").Value = 86.317567567567565
$ws.Range("A8").Value = 87.331081081081081
$ws.Range("A9").Value = 87.5
$ws.Range("A10").Value = 87.837837837837839
$ws.Range("A11").Value = 87.668918918918919
$ws.Range("A12").Value = 86.824324324324323
$ws.Range("A13").Value = 86.486486486486484
$ws.Range("A14").Value = 86.824324324324323
$ws.Range("A15").Value = 87.331081081081081
$ws.Range("A16").Value = 87.162162162162161
$ws.Range("A17").Value = 86.655405405405403
$ws.Range("A18").Value = 87.162162162162161
$ws.Range("A19").Value = 86.317567567567565
$ws.Range("A20").Value = 87.331081081081081
$ws.Range("A21").Value = 87.331081081081081
$ws.Range("A22").Value = 87.5
$ws.Range("A23").Value = 86.148648648648646
$ws.Range("A24").Value = 85.979729729729726
$ws.Range("A25").Value = 85.641891891891902
$ws.Range("A26").Value = 86.486486486486484
$ws.Range("A27").Value = 86.148648648648646
$ws.Range("A28").Value = 86.486486486486484
$ws.Range("A29").Value = 87.162162162162161
$ws.Range("A30").Value = 87.331081081081081
$ws.Range("A31").Value = 86.824324324324323
$ws.Range("A32").Value = 86.655405405405403
$ws.Range("A33").Value = 86.824324324324323
$ws.Range("A34").Value = 86.824324324324323
$ws.Range("A35").Value = 86.824324324324323
$ws.Range("A36").Value = 86.824324324324323
$ws.Range("A37").Value = 88.175675675675677
$ws.Range("A38").Value = 86.317567567567565
$ws.Range("A39").Value = 86.993243243243242
$ws.Range("A40").Value = 86.486486486486484
$ws.Range("A41").Value = 86.655405405405403
$ws.Range("A42").Value = 86.486486486486484
$ws.Range("A43").Value = 86.317567567567565
$ws.Range("A44").Value = 86.655405405405403
$ws.Range("A45").Value = 86.824324324324323
$ws.Range("A46").Value = 86.655405405405403
$ws.Range("A47").Value = 87.5
$ws.Range("A48").Value = 87.5
$ws.Range("A49").Value = 87.668918918918919
